$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1574.12
$ws.Range("I33").Value = 419.55
$ws.Range("J33").Value = 6192.4
$ws.Range("K33").Value = 419.55
$ws.Range("L33").Value = 6192.4
$ws.Range("M33").Value = -190.55
$ws.Range("N33").Value = -6650.4

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 1925
$ws.Range("I34").Value = 1925
$ws.Range("K34").Value = 1925
$ws.Range("M34").Value = -1654
$ws.Range("H74").Value = 5421.125
$ws.Range("I74").Value = 3559.9
$ws.Range("K74").Value = 3559.9
$ws.Range("M74").Value = -2685.9
$ws.Range("H77").Value = 5421.125
$ws.Range("I77").Value = 3559.9
$ws.Range("K77").Value = 17799.5
$ws.Range("M77").Value = -13431.5
$ws.Range("H122").Value = 3706.5334
$ws.Range("I122").Value = 3520.1667
$ws.Range("K122").Value = 10560.5001
$ws.Range("M122").Value = -8110.500100000001
$ws.Range("H132").Value = 1481.1522
$ws.Range("I132").Value = 1051.2972
$ws.Range("J132").Value = 3248.3333
$ws.Range("K132").Value = 3153.8916
$ws.Range("L132").Value = 9744.999899999999
$ws.Range("M132").Value = -623.8915999999999
$ws.Range("N132").Value = -14804.9999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 10376.267
$ws.Range("I20").Value = 2466.375
$ws.Range("J20").Value = 19416.143
$ws.Range("K20").Value = 2466.375
$ws.Range("L20").Value = 19416.143
$ws.Range("M20").Value = -2219.375
$ws.Range("N20").Value = -19910.143
$ws.Range("H80").Value = 40053.2
$ws.Range("J80").Value = 66.5
$ws.Range("L80").Value = 66.5
$ws.Range("N80").Value = -2062.5
$ws.Range("H83").Value = 40053.2
$ws.Range("J83").Value = 66.5
$ws.Range("L83").Value = 332.5
$ws.Range("N83").Value = -10316.5
$ws.Range("H86").Value = 2398.1428
$ws.Range("I86").Value = 1999.3334
$ws.Range("J86").Value = 2697.25
$ws.Range("K86").Value = 1999.3334
$ws.Range("L86").Value = 2697.25
$ws.Range("M86").Value = -876.3334
$ws.Range("N86").Value = -4943.25
$ws.Range("H89").Value = 2398.1428
$ws.Range("I89").Value = 1999.3334
$ws.Range("J89").Value = 2697.25
$ws.Range("K89").Value = 9996.666999999999
$ws.Range("L89").Value = 13486.25
$ws.Range("M89").Value = -4380.666999999999
$ws.Range("N89").Value = -24718.25
$ws.Range("H107").Value = 1710.6428
$ws.Range("J107").Value = 2368.2856
$ws.Range("L107").Value = 2368.2856
$ws.Range("N107").Value = -6208.2856
$ws.Range("H115").Value = 49999
$ws.Range("I115").Value = 49999
$ws.Range("K115").Value = 49999
$ws.Range("M115").Value = -48432

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 10499.5
$ws.Range("J43").Value = 10499.5
$ws.Range("L43").Value = 10499.5
$ws.Range("N43").Value = -10867.5
$ws.Range("H92").Value = 60300
$ws.Range("J92").Value = 60300
$ws.Range("L92").Value = 60300
$ws.Range("N92").Value = -65292
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("H101").Value = 10499.5
$ws.Range("J101").Value = 10499.5
$ws.Range("L101").Value = 10499.5
$ws.Range("N101").Value = -16989.5
$ws.Range("H105").Value = 1448.421
$ws.Range("I105").Value = 1644.3572
$ws.Range("K105").Value = 1644.3572
$ws.Range("M105").Value = 102.6428000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 252
$ws.Range("I13").Value = 252
$ws.Range("K13").Value = 756
$ws.Range("M13").Value = -588
$ws.Range("H17").Value = 787.125
$ws.Range("J17").Value = 933.3333
$ws.Range("L17").Value = 2799.9999
$ws.Range("N17").Value = -3137.9999
$ws.Range("H34").Value = 10259.857
$ws.Range("J34").Value = 11869.75
$ws.Range("L34").Value = 35609.25
$ws.Range("N34").Value = -35777.25
$ws.Range("H39").Value = 8056.4165
$ws.Range("J39").Value = 9519.666999999999
$ws.Range("L39").Value = 28559.001
$ws.Range("N39").Value = -29147.001
$ws.Range("H55").Value = 11549.3
$ws.Range("J55").Value = 11549.3
$ws.Range("L55").Value = 34647.89999999999
$ws.Range("N55").Value = -35001.89999999999
$ws.Range("H102").Value = 2999.5
$ws.Range("I102").Value = 1000
$ws.Range("J102").Value = 4999
$ws.Range("K102").Value = 3000
$ws.Range("L102").Value = 14997
$ws.Range("M102").Value = -566
$ws.Range("N102").Value = -19865
$ws.Range("H137").Value = 5281.684
$ws.Range("I137").Value = 2077
$ws.Range("J137").Value = 8842.444
$ws.Range("K137").Value = 6231
$ws.Range("L137").Value = 26527.332
$ws.Range("M137").Value = -1131
$ws.Range("N137").Value = -36727.33199999999
$ws.Range("H138").Value = 3140.0625
$ws.Range("I138").Value = 2247.125
$ws.Range("K138").Value = 6741.375
$ws.Range("M138").Value = -1601.375

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 838
$ws.Range("I2").Value = 900.8095
$ws.Range("K2").Value = 900.8095
$ws.Range("M2").Value = -787.8095
$ws.Range("H32").Value = 54797.145
$ws.Range("J32").Value = 56430
$ws.Range("L32").Value = 56430
$ws.Range("N32").Value = -57022
$ws.Range("H42").Value = 42145
$ws.Range("J42").Value = 42145
$ws.Range("L42").Value = 42145
$ws.Range("N42").Value = -43115
$ws.Range("H45").Value = 89130.39999999999
$ws.Range("J45").Value = 92663
$ws.Range("L45").Value = 92663
$ws.Range("N45").Value = -93781
$ws.Range("H51").Value = 94331.5
$ws.Range("J51").Value = 100775.336
$ws.Range("L51").Value = 100775.336
$ws.Range("N51").Value = -101793.336
$ws.Range("H102").Value = 2466.2173
$ws.Range("I102").Value = 1401.2667
$ws.Range("K102").Value = 1401.2667
$ws.Range("M102").Value = 220.7333000000001
$ws.Range("H103").Value = 49995
$ws.Range("J103").Value = 49995
$ws.Range("L103").Value = 49995
$ws.Range("N103").Value = -52339
$ws.Range("H115").Value = 42145
$ws.Range("J115").Value = 42145
$ws.Range("L115").Value = 42145
$ws.Range("N115").Value = -44495
$ws.Range("H122").Value = 5065.44
$ws.Range("I122").Value = 2481.0833
$ws.Range("K122").Value = 7443.249899999999
$ws.Range("M122").Value = -4993.249899999999
$ws.Range("H123").Value = 54775.223
$ws.Range("J123").Value = 54775.223
$ws.Range("L123").Value = 54775.223
$ws.Range("N123").Value = -59675.223
$ws.Range("H126").Value = 7404.6
$ws.Range("I126").Value = 7409.5713
$ws.Range("K126").Value = 22228.7139
$ws.Range("M126").Value = -19758.7139
$ws.Range("H132").Value = 4026.6382
$ws.Range("I132").Value = 2401.3845
$ws.Range("K132").Value = 7204.1535
$ws.Range("M132").Value = -4674.1535

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 16999.666
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H132").Value = 6388.477
$ws.Range("I132").Value = 4693.1665
$ws.Range("K132").Value = 14079.4995
$ws.Range("M132").Value = -11549.4995
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 3214
$ws.Range("I136").Value = 2927.9062
$ws.Range("J136").Value = 4046.2727
$ws.Range("K136").Value = 8783.7186
$ws.Range("L136").Value = 12138.8181
$ws.Range("M136").Value = -6233.7186
$ws.Range("N136").Value = -17238.8181
